$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E7").Value = "Valid account data"
$ws.Range("F7").Value = '("1000", "C001", 1000, date(2010, 1, 1), 2.00)'
$ws.Range("G7").Value = 'account_number="1000", balance=1000, date_created=date(2010, 1, 1), management_fee=2.00'

$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = '("1000", "C001", 1000, date(2015, 6, 1), "invalid")'
$ws.Range("G8").Value = "management_fee defaults to 2.55"

$ws.Range("E9").Value = "self.date_created=date(2000, 1, 1)"
$ws.Range("F9").Value = "get_service_charges()"
$ws.Range("G9").Value = "Returns base service charge of 0.50 (management fee waived)"

$ws.Range("E10").Value = "self.date_created=date(2015, 3, 16)"
$ws.Range("F10").Value = "get_service_charges()"
$ws.Range("G10").Value = "Returns 2.50 (0.50 + 2.00, assuming management fee applies)"

$ws.Range("E11").Value = "self.date_created=date(2020, 1, 1)"
$ws.Range("F11").Value = "get_service_charges()"
$ws.Range("G11").Value = "Returns 2.55 (0.50 + 2.55, assuming default management fee)"

$ws.Range("E12").Value = "self.date_created=date(2000, 1, 1)"
$ws.Range("F12").Value = "str(investment_account)"
$ws.Range("G12").Value = "Returns ""Account Number: 1000 Balance: `$1,000.00`nDate Created: 2000-01-01 Management Fee: Waived Account Type: Inve"""

$ws.Range("E13").Value = "self.date_created=date(2020, 1, 1)"
$ws.Range("F13").Value = "str(investment_account)"
$ws.Range("G13").Value = "Returns ""Account Number: 1000 Balance: `$1,000.00`nDate Created: 2020-01-01 Management Fee: `$2.55 Account Type: Inve"""

$ws.Range("G13").Select()
